$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52-173 down to 53-174
$ws.Rows(52).Insert()

# Fill the new row 52 with data. Columns A,B,C,E,F,G,H,Q,R are copied down automatically
# by the insert (since Excel's Insert shifts the old row 52 content down to row 53, leaving
# the new row 52 blank), so we set them explicitly to match the required final values.
$ws.Range("A52").Value = 8
$ws.Range("B52").Value = "Terminal La Palmera de La Serena"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 45259
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 100112028
$ws.Range("G52").Value = "Sandia"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 700
$ws.Range("L52").Value = 800
$ws.Range("M52").Value = 750
$ws.Range("N52").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O52").Value = "Perú"
$ws.Range("P52").Value = 750
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
